$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 48, shifting existing rows 48..125 down to 49..126.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with its data (same template as the
# surrounding Ajo/Chino/Primera rows, with its own Fecha/Volumen/Precio values).
$ws.Range("A48").Value = 7
$ws.Range("B48").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C48").Value = "Ñuble"
$ws.Range("D48").Value = 44469
$ws.Range("E48").Value = 16
$ws.Range("F48").Value = 100112003
$ws.Range("G48").Value = "Ajo"
$ws.Range("H48").Value = "Chino"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 60
$ws.Range("K48").Value = 16000
$ws.Range("L48").Value = 17000
$ws.Range("M48").Value = 16500
$ws.Range("N48").Value = "$/caja 10 kilos"
$ws.Range("O48").Value = "China"
$ws.Range("P48").Value = 1650
$ws.Range("Q48").Value = 10
$ws.Range("R48").Value = "Hortaliza"
